# Rename the diff-report header columns to reference the actual
# format-version names (FV2404 / FV2410) instead of the generic
# "_old" / "_new" suffixes, then turn the sheet into a proper Excel
# Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J: "<name>_old" -> "<name>_FV2404"
for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseHeaders[$i])_FV2404"
}

# Column K stays "diff"

# Columns L-U: "<name>_new" -> "<name>_FV2410"
for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseHeaders[$i])_FV2410"
}

# Turn the used range into an Excel Table ("Table1") with the
# renamed headers, keeping the existing default table style.
$tableRange = $ws.Range("A1:U93")
$lo = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1, keep top-left at A2).
$ws.Activate()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
